$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    12 = 1
    15 = 6
    16 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 2
    21 = 2
    25 = 6
    28 = 4
    29 = 2
    36 = 2
    37 = 2
    38 = 2
    42 = 6
    45 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("Q$row").Value = $values[$row]
}
